# StudyingHours.xlsx — "Add files via upload" edit
#
# The commit re-saved the workbook from a newer Excel build (refreshed
# fileVersion/calcPr/theme metadata, a new empty xl/persons/person.xml,
# refreshed window geometry, etc.) and, as the actual content change,
# filled in the previously-empty Avg_Hours cell for row 5 ("Survey 3",
# column B = B5) with 7.53 and left the selection sitting on that cell.
#
# Only the latter is reachable through the Excel object model — the rest
# is save-time/application-version plumbing Excel regenerates on its own
# and isn't exposed as settable COM properties — so we reproduce the
# actual data edit here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Survey 3") previously had no Avg_Hours (column B) value.
$ws.Range("B5").Value = 7.53

# Move/leave the active selection on the cell that was just edited.
$ws.Range("B5").Select() | Out-Null
